# Atualização de bases das ligas, do dia: 14-06-2024 às 20:31
# Swap the full data (columns B..AD) between each of the following row pairs.
# Column A (the running index) stays put at its row position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$pairs = @(
    @(86, 87),
    @(147, 148),
    @(221, 222)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Capture both rows' values first so the swap doesn't clobber itself.
    $row1Vals = @{}
    $row2Vals = @{}
    foreach ($col in $cols) {
        $row1Vals[$col] = $ws.Range("$col$r1").Value2
        $row2Vals[$col] = $ws.Range("$col$r2").Value2
    }

    foreach ($col in $cols) {
        $ws.Range("$col$r1").Value2 = $row2Vals[$col]
        $ws.Range("$col$r2").Value2 = $row1Vals[$col]
    }
}
